# "add monte_carlo and update database"
# Populate the yearly income-statement figures (rial) for Sheranol/Ravankar
# on the "Overview" sheet. Columns D:H correspond to the five fiscal
# periods (1396/12 .. 1400/12); rows 11-27 are the statement line items,
# previously all zero / placeholder "-" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# فروش (Sales)
$ws.Range("D11").Value = 12944707
$ws.Range("E11").Value = 20919381
$ws.Range("F11").Value = 30623002
$ws.Range("G11").Value = 51595247
$ws.Range("H11").Value = 87924869

# بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -9072135
$ws.Range("E12").Value = -16505656
$ws.Range("F12").Value = -25341376
$ws.Range("G12").Value = -36534778
$ws.Range("H12").Value = -68089849

# سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 3872572
$ws.Range("E13").Value = 4413725
$ws.Range("F13").Value = 5281626
$ws.Range("G13").Value = 15060469
$ws.Range("H13").Value = 19835020

# هزینه های عمومی, اداری و تشکیلاتی (General & administrative expenses)
$ws.Range("D14").Value = -949290
$ws.Range("E14").Value = -1402872
$ws.Range("F14").Value = -1457654
$ws.Range("G14").Value = -1889165
$ws.Range("H14").Value = -2385117

# هزینه کاهش ارزش دریافتنی‌ها (Impairment expense) - D15 stays "-"
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net)
$ws.Range("D16").Value = -26975
$ws.Range("E16").Value = 567763
$ws.Range("F16").Value = 83586
$ws.Range("G16").Value = -855868
$ws.Range("H16").Value = 679363

# سود (زیان) عملیاتی (Operating profit)
$ws.Range("D17").Value = 2896307
$ws.Range("E17").Value = 3578616
$ws.Range("F17").Value = 3907558
$ws.Range("G17").Value = 12315436
$ws.Range("H17").Value = 18129266

# هزینه های مالی (Finance costs)
$ws.Range("D18").Value = -390262
$ws.Range("E18").Value = -476394
$ws.Range("F18").Value = -847533
$ws.Range("G18").Value = -1229126
$ws.Range("H18").Value = -1870249

# خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net)
$ws.Range("D19").Value = 12785
$ws.Range("E19").Value = 154516
$ws.Range("F19").Value = 316920
$ws.Range("G19").Value = 553714
$ws.Range("H19").Value = 485700

# سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing operations)
$ws.Range("D20").Value = 2518830
$ws.Range("E20").Value = 3256738
$ws.Range("F20").Value = 3376945
$ws.Range("G20").Value = 11640024
$ws.Range("H20").Value = 16744717

# مالیات (Tax)
$ws.Range("D21").Value = -397567
$ws.Range("E21").Value = -478089
$ws.Range("F21").Value = -354451
$ws.Range("G21").Value = -950685
$ws.Range("H21").Value = -2181134

# سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing operations)
$ws.Range("D22").Value = 2121263
$ws.Range("E22").Value = 2778649
$ws.Range("F22").Value = 3022494
$ws.Range("G22").Value = 10689339
$ws.Range("H22").Value = 14563583

# سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (Discontinued operations, net of tax)
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# سود (زیان) خالص (Net profit)
$ws.Range("D24").Value = 2121263
$ws.Range("E24").Value = 2778649
$ws.Range("F24").Value = 3022494
$ws.Range("G24").Value = 10689339
$ws.Range("H24").Value = 14563583

# سود هر سهم پس از کسر مالیات (EPS after tax) - E25/F25 relationship: E25 stays "-"
$ws.Range("D25").Value = 1061
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 5345
$ws.Range("H25").Value = 7282

# سرمایه (Capital) - E26 becomes numeric 0
$ws.Range("D26").Value = 2000000
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2000000
$ws.Range("G26").Value = 2000000
$ws.Range("H26").Value = 2000000

# سود هر سهم بر اساس آخرین سرمایه (EPS based on latest capital)
$ws.Range("D27").Value = 1061
$ws.Range("E27").Value = 1389
$ws.Range("F27").Value = 1511
$ws.Range("G27").Value = 5345
$ws.Range("H27").Value = 7282
